$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.147.93"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.053.97"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.52"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.01"
$ws.Range("E7").Value = "  -5.79%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0787"
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.108"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.77"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "2.354.59"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.834"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.79"
$ws.Range("E15").Value = "  +6.69%  "
$ws.Range("D16").Value = "2.050.46"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.09"
$ws.Range("E17").Value = "  +19.16%  "
$ws.Range("D18").Value = "37.179.18"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.09"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.41"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.96"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +3.19%  "
$ws.Range("E25").Value = "  +5.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.37"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.15"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0907"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.75"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.13"
$ws.Range("E40").Value = "  +11.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.11"
$ws.Range("E41").Value = "  +12.09%  "
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.16"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.38"
$ws.Range("E44").Value = "  -7.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.69"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.48"
$ws.Range("E46").Value = "  -3.82%  "
$ws.Range("D47").Value = "1.289.50"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.90"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "2.246.70"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.63"
$ws.Range("E51").Value = "  -18.60%  "
